$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '27.418.69'
$ws.Range('E2').Value = '  +0.25%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.635.65'
$ws.Range('E3').Value = '  -0.75%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '212.30'
$ws.Range('E5').Value = '  -0.70%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.530'
$ws.Range('E6').Value = '  +4.28%  '
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '22.93'
$ws.Range('E8').Value = '  -4.16%  '
$ws.Range('E9').Value = '  -2.15%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0609'
$ws.Range('E10').Value = '  -0.82%  '
$ws.Range('E11').Value = '  +1.08%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.866.98'
$ws.Range('E12').Value = '  -0.88%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.624.10'
$ws.Range('E13').Value = '  -1.60%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.578'
$ws.Range('E14').Value = '  +3.04%  '
$ws.Range('E15').Value = '  -1.91%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '64.05'
$ws.Range('E16').Value = '  -2.25%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '27.415.74'
$ws.Range('E17').Value = '  +0.16%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '228.90'
$ws.Range('E18').Value = '  -2.32%  '
$ws.Range('E19').Value = '  -0.42%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.54'
$ws.Range('E20').Value = '  +0.74%  '
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.30'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.72'
$ws.Range('E23').Value = '  +5.63%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.95'
$ws.Range('E24').Value = '  -3.74%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '149.36'
$ws.Range('E25').Value = '  +2.48%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '6.99'
$ws.Range('E26').Value = '  -2.22%  '
$ws.Range('E27').Value = '  +1.78%  '
$ws.Range('E28').Value = '  -0.10%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.53'
$ws.Range('E29').Value = '  -3.06%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.19'
$ws.Range('E30').Value = '  -0.40%  '
$ws.Range('E31').Value = '  -1.92%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.28'
$ws.Range('E32').Value = '  -0.25%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.17'
$ws.Range('E33').Value = '  +3.40%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.410.76'
$ws.Range('E34').Value = '  -2.83%  '
$ws.Range('E35').Value = '  +2.19%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.35'
$ws.Range('E36').Value = '  -2.00%  '
$ws.Range('E37').Value = '  -0.19%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.873'
$ws.Range('E38').Value = '  -3.95%  '
$ws.Range('E39').Value = '  -1.77%  '
$ws.Range('E40').Value = '  -1.41%  '
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.832'
$ws.Range('E42').Value = '  +5.79%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.50'
$ws.Range('E43').Value = '  +1.14%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.23'
$ws.Range('E44').Value = '  +0.36%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '64.68'
$ws.Range('E45').Value = '  -0.73%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.776.45'
$ws.Range('E46').Value = '  -0.79%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.65'
$ws.Range('E47').Value = '  -3.20%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '85.13'
$ws.Range('E48').Value = '  -3.54%  '
$ws.Range('E49').Value = '  -0.16%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0989'
$ws.Range('E50').Value = '  -1.64%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.71'
$ws.Range('E51').Value = '  -0.90%  '
